# Add a new data row (row 4) to the EMP_RECORD sheet, matching the
# existing sheet's convention of storing every field - including the
# numeric-looking ID/serial columns A, B, I - as text rather than numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMP_RECORD")

$rowValues = @{
    1 = "3"
    2 = "23"
    3 = "asdfl;k"
    4 = "saflk;j"
    5 = "salkfjlk"
    6 = "slkfaj"
    7 = "lksfdajlk"
    8 = "lkdsafj"
    9 = "33"
}

# Columns whose values look numeric ("3", "23", "33") need the cell
# pre-formatted as Text, otherwise Excel auto-coerces the assignment into
# a real number (as it already does for row 2 in this sheet). The other
# columns are plain words, so a normal value assignment already stores
# them as text - no need to touch their formatting.
$textForcedColumns = @(1, 2, 9)

foreach ($col in 1..9) {
    $cell = $ws.Cells.Item(4, $col)
    if ($textForcedColumns -contains $col) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $rowValues[$col]
}
